$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the stored credential e-mail with the new one and turn it into a
# live mailto hyperlink (Excel auto-applies the built-in "Hyperlink" cell
# style, adds the Hyperlink font/cellStyle to styles.xml and records the
# link in the worksheet's relationships / <hyperlinks> part).
$ws.Range("A2").Hyperlinks.Add($ws.Range("A2"), "mailto:julieth27f@hotmail.com", "", "", "julieth27f@hotmail.com") | Out-Null

# Replace the numeric password placeholder with the real password text.
$ws.Range("B2").Value = "pi79810859nk"

# Move the active selection like the recorded session (cosmetic, matches
# the saved sheetView/selection in the target file).
$ws.Range("G8").Select() | Out-Null

Write-Output "done"
